$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "Anselmo-Gestao Intregrada"

$ws.Range("B4").Value = "[-, -, -, Aline S. M.-T. M. Metalicos-1B]"
$ws.Range("E4").Value = "Anselmo-Gestao Intregrada"
$ws.Range("F4").Value = "[Emerson-Comandos Eletricos-1B, Emerson-Comandos Eletricos-1B, Emerson-Comandos Eletricos-1B, Emerson-Comandos Eletricos-1B]"

$ws.Range("B6").Value = "[Ernane-Desenho tecnico mecanico-1B, Ernane-Desenho tecnico mecanico-1B, Ernane-Desenho tecnico mecanico-1B, Aline S. M.-T. M. Metalicos-1B]"
$ws.Range("C6").Value = "[Weslei-Metrologia 1-1B, Weslei-Metrologia 1-1B, Weslei-Metrologia 1-1B, Weslei-Metrologia 1-1B]"
